# Update "time_taken" (column F) values on the "data" sheet with the
# freshly re-queried timestamps from the refreshed PanelApp API pull.
$newTimeTaken = @(
    "2021-10-05 14:33:32.647244",
    "2021-10-05 14:33:32.647252",
    "2021-10-05 14:33:32.647255",
    "2021-10-05 14:33:32.647258",
    "2021-10-05 14:33:32.647260",
    "2021-10-05 14:33:32.647263",
    "2021-10-05 14:33:32.647265",
    "2021-10-05 14:33:32.647268",
    "2021-10-05 14:33:32.647270",
    "2021-10-05 14:33:32.647273",
    "2021-10-05 14:33:32.647275",
    "2021-10-05 14:33:32.647277",
    "2021-10-05 14:33:32.647280",
    "2021-10-05 14:33:32.647282",
    "2021-10-05 14:33:32.647285",
    "2021-10-05 14:33:32.647287",
    "2021-10-05 14:33:32.647290",
    "2021-10-05 14:33:32.647292",
    "2021-10-05 14:33:32.647295",
    "2021-10-05 14:33:32.647297",
    "2021-10-05 14:33:32.647300",
    "2021-10-05 14:33:32.647302",
    "2021-10-05 14:33:32.647305",
    "2021-10-05 14:33:32.647307",
    "2021-10-05 14:33:32.647310",
    "2021-10-05 14:33:32.647312",
    "2021-10-05 14:33:32.647314",
    "2021-10-05 14:33:32.647317",
    "2021-10-05 14:33:32.647319",
    "2021-10-05 14:33:32.647322",
    "2021-10-05 14:33:32.647324",
    "2021-10-05 14:33:32.647326",
    "2021-10-05 14:33:32.647329",
    "2021-10-05 14:33:32.647332",
    "2021-10-05 14:33:32.647334",
    "2021-10-05 14:33:32.647336",
    "2021-10-05 14:33:32.647339",
    "2021-10-05 14:33:32.647341",
    "2021-10-05 14:33:32.647344",
    "2021-10-05 14:33:32.647346",
    "2021-10-05 14:33:32.647349",
    "2021-10-05 14:33:32.647352",
    "2021-10-05 14:33:32.647354",
    "2021-10-05 14:33:32.647356",
    "2021-10-05 14:33:32.647359",
    "2021-10-05 14:33:32.647362",
    "2021-10-05 14:33:32.647364",
    "2021-10-05 14:33:32.647366",
    "2021-10-05 14:33:32.647369",
    "2021-10-05 14:33:32.647371",
    "2021-10-05 14:33:32.647374",
    "2021-10-05 14:33:32.647376",
    "2021-10-05 14:33:32.647379",
    "2021-10-05 14:33:32.647382",
    "2021-10-05 14:33:32.647384",
    "2021-10-05 14:33:32.647387",
    "2021-10-05 14:33:32.647389",
    "2021-10-05 14:33:32.647391",
    "2021-10-05 14:33:32.647394",
    "2021-10-05 14:33:32.647396",
    "2021-10-05 14:33:32.647399",
    "2021-10-05 14:33:32.647401",
    "2021-10-05 14:33:32.647404",
    "2021-10-05 14:33:32.647406",
    "2021-10-05 14:33:32.647410",
    "2021-10-05 14:33:32.647413",
    "2021-10-05 14:33:32.647415",
    "2021-10-05 14:33:32.647418",
    "2021-10-05 14:33:32.647420",
    "2021-10-05 14:33:32.647422",
    "2021-10-05 14:33:32.647425",
    "2021-10-05 14:33:32.647428",
    "2021-10-05 14:33:32.647430",
    "2021-10-05 14:33:32.647432",
    "2021-10-05 14:33:32.647435",
    "2021-10-05 14:33:32.647438",
    "2021-10-05 14:33:32.647442",
    "2021-10-05 14:33:32.647446",
    "2021-10-05 14:33:32.647448",
    "2021-10-05 14:33:32.647450",
    "2021-10-05 14:33:32.647453",
    "2021-10-05 14:33:32.647456",
    "2021-10-05 14:33:32.647458",
    "2021-10-05 14:33:32.647461",
    "2021-10-05 14:33:32.647463",
    "2021-10-05 14:33:32.647466",
    "2021-10-05 14:33:32.647469",
    "2021-10-05 14:33:32.647471",
    "2021-10-05 14:33:32.647474",
    "2021-10-05 14:33:32.647476",
    "2021-10-05 14:33:32.647478",
    "2021-10-05 14:33:32.647481",
    "2021-10-05 14:33:32.647485",
    "2021-10-05 14:33:32.647488",
    "2021-10-05 14:33:32.647490",
    "2021-10-05 14:33:32.647493",
    "2021-10-05 14:33:32.647495",
    "2021-10-05 14:33:32.647498",
    "2021-10-05 14:33:32.647500",
    "2021-10-05 14:33:32.647503",
    "2021-10-05 14:33:32.647505",
    "2021-10-05 14:33:32.647508",
    "2021-10-05 14:33:32.647510"
)

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

for ($row = 2; $row -le 104; $row++) {
    $dataSheet.Cells.Item($row, 6).Value = $newTimeTaken[$row - 2]
}

# Add a new "metadata" tab (right after "data") describing the panel pull
# itself: which PanelApp panel/version was queried, when, and via which
# request URL.
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Re-use the existing bold/bordered/centered header style (style index 1,
# already present in the workbook) instead of minting a brand-new style --
# copy formatting from the "data" sheet's own header cell.
$headerStyleSource = $dataSheet.Range("B1")
$headerStyleSource.Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$headerStyleSource.Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Congenital abnormalities of the kidneys and urinary tract (CAKUT)_SuperPanel"
$metaSheet.Range("C2").Value = 251

# Force D2 to be stored as literal text "0.226" (not coerced to the number
# 0.226) to match the source panel metadata's string representation, then
# drop the temporary text format so the cell keeps the default style.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.226"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-10-04T06:48:28.556180Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:32.643944"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/251/?format=json"

# Keep "data" as the active/selected sheet (unchanged bookViews/activeTab),
# since only a new tab was appended -- the view focus itself didn't move.
$dataSheet.Activate()
